# python/xlsx/python XL test.xlsx -- refresh the provenance sheet's
# machine/environment info (directory, python version, $USER, $HOSTNAME,
# $HOME) and the run timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# directory
$ws.Range("B5").Value = "/Volumes/Tlaltecuhtli/repos/GitHub/topa-development/python/xlsx"

# python version (two lines, joined with an embedded newline)
$ws.Range("B6").Value = "3.7.0 (default, Jun 28 2018, 07:39:16) `n[Clang 4.0.1 (tags/RELEASE_401/final)]"

# $USER
$ws.Range("B9").Value = "l127914"

# $HOSTNAME
$ws.Range("B10").Value = "Cauchy.Schwarz"

# $HOME
$ws.Range("B11").Value = "/Users/l127914"

# timestamp (serial date/time number, cell already carries the date format)
$ws.Range("B12").Value = 43438.43873602271
